$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-13 22:48:27'
$ws.Range("I2").Value = '4.5 mm'
$ws.Range("E3").Value = '2026-02-13 22:48:30'
$ws.Range("I3").Value = '7.8 mm'
$ws.Range("E4").Value = '2026-02-13 22:48:33'
$ws.Range("J4").Value = '993.3 hPa'
$ws.Range("L4").Value = '26.3 km/h - 273º 22:08 TU'
$ws.Range("E5").Value = '2026-02-13 22:48:35'
$ws.Range("I5").Value = '3.8 mm'
$ws.Range("N5").Value = '-4.2 °C 22:29 TU'
$ws.Range("E6").Value = '2026-02-13 22:48:38'
$ws.Range("J6").Value = '993.3 hPa'
$ws.Range("E7").Value = '2026-02-13 22:48:40'
$ws.Range("J7").Value = '993.7 hPa'
$ws.Range("E8").Value = '2026-02-13 22:48:43'
$ws.Range("J8").Value = '993.6 hPa'
$ws.Range("E9").Value = '2026-02-13 22:48:45'
$ws.Range("E10").Value = '2026-02-13 22:48:48'
$ws.Range("H10").Value = "'89%"
$ws.Range("E11").Value = '2026-02-13 22:48:50'
$ws.Range("E12").Value = '2026-02-13 22:48:53'
$ws.Range("O12").Value = '9.6 °C'
$ws.Range("E13").Value = '2026-02-13 22:48:55'
$ws.Range("E14").Value = '2026-02-13 22:48:58'
$ws.Range("H14").Value = "'83%"
$ws.Range("L14").Value = '48.2 km/h - 287º 22:01 TU'
$ws.Range("E15").Value = '2026-02-13 22:49:00'
$ws.Range("I15").Value = '5.6 mm'
$ws.Range("E16").Value = '2026-02-13 22:49:03'
$ws.Range("I16").Value = '14.7 mm'
$ws.Range("E17").Value = '2026-02-13 22:49:05'
$ws.Range("E18").Value = '2026-02-13 22:49:08'
$ws.Range("J18").Value = '993.5 hPa'
$ws.Range("E19").Value = '2026-02-13 22:49:10'
$ws.Range("H19").Value = "'91%"
$ws.Range("O19").Value = '3.7 °C'
$ws.Range("E20").Value = '2026-02-13 22:49:13'
$ws.Range("I20").Value = '24.4 mm'
$ws.Range("E21").Value = '2026-02-13 22:49:16'
$ws.Range("J21").Value = '996.5 hPa'
$ws.Range("N21").Value = '-0.3 °C 22:17 TU'
$ws.Range("E22").Value = '2026-02-13 22:49:18'
$ws.Range("L22").Value = '65.5 km/h - 301º 22:27 TU'
$ws.Range("E23").Value = '2026-02-13 22:49:21'
$ws.Range("G23").Value = '192 cm'
$ws.Range("I23").Value = '13.8 mm'
$ws.Range("E24").Value = '2026-02-13 22:49:23'
$ws.Range("J24").Value = '994.7 hPa'
$ws.Range("E25").Value = '2026-02-13 22:49:26'
$ws.Range("I25").Value = '10.2 mm'
$ws.Range("E26").Value = '2026-02-13 22:49:29'
$ws.Range("E27").Value = '2026-02-13 22:49:31'
$ws.Range("E28").Value = '2026-02-13 22:49:34'
$ws.Range("H28").Value = "'82%"
$ws.Range("J28").Value = '993.8 hPa'
$ws.Range("O28").Value = '6.4 °C'
$ws.Range("E29").Value = '2026-02-13 22:49:36'
$ws.Range("E30").Value = '2026-02-13 22:49:39'
$ws.Range("H30").Value = "'77%"
$ws.Range("J30").Value = '993.2 hPa'
$ws.Range("E31").Value = '2026-02-13 22:49:41'
$ws.Range("I31").Value = '5.0 mm'
$ws.Range("J31").Value = '992.2 hPa'
$ws.Range("E32").Value = '2026-02-13 22:49:44'
$ws.Range("E33").Value = '2026-02-13 22:49:47'
$ws.Range("H33").Value = "'91%"
$ws.Range("I33").Value = '5.0 mm'
$ws.Range("J33").Value = '995.3 hPa'
$ws.Range("E34").Value = '2026-02-13 22:49:49'
$ws.Range("E35").Value = '2026-02-13 22:49:52'
$ws.Range("H35").Value = "'78%"
$ws.Range("J35").Value = '994.9 hPa'
$ws.Range("L35").Value = '86.0 km/h - 265º 22:14 TU'
$ws.Range("E36").Value = '2026-02-13 22:49:54'
$ws.Range("J36").Value = '993.3 hPa'
$ws.Range("O36").Value = '10.8 °C'
$ws.Range("E37").Value = '2026-02-13 22:49:57'
$ws.Range("J37").Value = '995.2 hPa'
$ws.Range("E38").Value = '2026-02-13 22:49:59'
$ws.Range("N38").Value = '7.7 °C 22:23 TU'
$ws.Range("E39").Value = '2026-02-13 22:50:02'
$ws.Range("I39").Value = '19.9 mm'
$ws.Range("E40").Value = '2026-02-13 22:50:04'
$ws.Range("J40").Value = '997.0 hPa'
$ws.Range("E41").Value = '2026-02-13 22:50:07'
$ws.Range("H41").Value = "'74%"
$ws.Range("J41").Value = '994.1 hPa'
$ws.Range("E42").Value = '2026-02-13 22:50:09'
$ws.Range("E43").Value = '2026-02-13 22:50:12'
$ws.Range("O43").Value = '6.2 °C'
$ws.Range("E44").Value = '2026-02-13 22:50:14'
$ws.Range("I44").Value = '11.0 mm'
$ws.Range("E45").Value = '2026-02-13 22:50:17'
$ws.Range("H45").Value = "'68%"
$ws.Range("O45").Value = '5.4 °C'
$ws.Range("E46").Value = '2026-02-13 22:50:19'
$ws.Range("J46").Value = '994.9 hPa'
$ws.Range("O46").Value = '9.3 °C'
